# The Google Sheet export regenerated the .ttl, which removed the three
# section-header rows that used to separate the blocks of the vocabulary
# template sheet ("Prefixes for controlled vocabularies..." at row 2,
# "Metadata about vocabulary" at row 7 and "Definition of terms" at row 15).
# Deleting those rows shifts everything below them up, which also drops the
# three now-unused blank "nicest-2-subjects:" template rows off the bottom
# of the used range (dimension shrinks from A1:S52 to A1:S49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows("15:15").Delete()
$ws.Rows("7:7").Delete()
$ws.Rows("2:2").Delete()
